$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rename the existing filepath headers to the "_en" variants.
$ws1.Range("G1").Value = "icon_filepath_en"
$ws1.Range("H1").Value = "audio_filepath_en"

# Insert two new columns (I, J) for the French-localized media columns;
# this pushes the existing "unique_id" column (and its data) from I to K.
$ws1.Columns("I:J").Insert()
$ws1.Range("I1").Value = "icon_filepath_fra"
$ws1.Range("J1").Value = "audio_filepath_fra"

# New header cells get a distinct (explicit black) font color/style.
$ws1.Range("I1:J1").Font.Color = 0

# Make the Modules_and_forms sheet the active tab, with I1:J1 selected.
$ws1.Activate()
$ws1.Range("I1:J1").Select()
